# Update the build/version string throughout the workbook.
$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# "About" sheet
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Sima Coal Mine, China, M1203, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: column S (build_version) for rows 2-8
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 8; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
